$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.353.22'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '1.711.63'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5303'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.06698'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2668'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.90'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07683'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('E12').Value = '  -2.12%  '
$ws.Range('D13').Value = '1.947.11'
$ws.Range('E13').Value = '  -0.63%  '
$ws.Range('D14').Value = '1.698.97'
$ws.Range('E14').Value = '  -1.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5839'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '0.0₅8236'
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.13'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '27.362.21'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '223.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.004'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.641'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.71%  '
$ws.Range('E22').Value = '  -2.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.67%  '
$ws.Range('E26').Value = '  -2.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1215'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.255'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.27'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05375'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.294'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.440'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.458'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.866'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9532'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.393'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5869'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01639'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.92%  '
$ws.Range('D40').Value = '1.091.39'
$ws.Range('E40').Value = '  +3.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.824'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.66%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.005'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8443'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.90'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('D45').Value = '1.854.36'
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('D46').Value = '0.0₈114'
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.96'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4533'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.007'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.117'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05234'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.34%  '
